$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the new, longer question titles
$ws.Columns.Item(1).ColumnWidth = 19.42578125

# Row 50: 560. Subarray Sum Equals K
$ws.Cells.Item(50, 1).Value = 560
$ws.Cells.Item(50, 2).Value = "Subarray Sum Equals K"
$ws.Cells.Item(50, 3).Value = "Java"
$ws.Cells.Item(50, 4).Value = 45049

# Row 51: Subarrays with XOR 'K' (Coding Ninja)
$ws.Cells.Item(51, 1).Value = "Coding Ninja"
$ws.Cells.Item(51, 2).Value = "Subarrays with XOR 'K'"
$ws.Cells.Item(51, 3).Value = "Java"
$ws.Cells.Item(51, 4).Value = 45049

# Match formatting of the row above (A column centered, D column date format)
$ws.Range("A50:D50").NumberFormat = $ws.Range("A49:D49").NumberFormat
$ws.Range("A51:D51").NumberFormat = $ws.Range("A49:D49").NumberFormat
$ws.Cells.Item(50, 1).HorizontalAlignment = -4108
$ws.Cells.Item(51, 1).HorizontalAlignment = -4108
$ws.Cells.Item(50, 4).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(51, 4).NumberFormat = "m/d/yyyy"

$ws.Range("E58").Select()
